# Regenerate merged AHB files
# - Rename the "_old" / "_new" header suffixes to the concrete version
#   identifiers "_FV2410" / "_FV2504" used by the merged AHB export.
# - Wrap the data range in an Excel Table (ListObject) with an AutoFilter.
# - Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row (row 1) -----------------------------------------
# Columns A:J carry the "_old" suffix -> "_FV2410"
for ($c = 1; $c -le 10; $c++) {
    $header = $ws.Cells.Item(1, $c).Value()
    $ws.Cells.Item(1, $c).Value = ($header -replace "_old$", "_FV2410")
}

# Column K is the static "diff" column header; columns L:U carry the
# "_new" suffix -> "_FV2504"
for ($c = 12; $c -le 21; $c++) {
    $header = $ws.Cells.Item(1, $c).Value()
    $ws.Cells.Item(1, $c).Value = ($header -replace "_new$", "_FV2504")
}

# --- 2) Turn the used range into an Excel Table ----------------------------
$dataRange = $ws.Range("A1:U63")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3) Freeze the header row ----------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
